{"js": "// Add \"<<Include>> Login, dopo cui \" before the actor reference in the\n// first step of each of the three use-case tables, turning:\n//   \"1. L'operatore sceglie ...\"      -> \"1. <<Include>> Login, dopo cui l'operatore sceglie ...\"\n//   \"1. L'utente sceglie ...\"         -> \"1. <<Include>> Login, dopo cui l'utente sceglie ...\"\n// (only the leading capital \"L\" of the actor word is replaced, turning it\n// lower-case and prefixing the \"<<Include>> Login, dopo cui \" phrase - the\n// rest of the sentence is untouched).\n\nconst body = context.document.body;\n\n// Each entry identifies the first step of a use case uniquely via the text\n// that immediately follows the actor reference we need to touch.\nconst targets = [\n  \"L\\u2019operatore sceglie la funzione \\u201cAggiungi\",\n  \"L\\u2019utente sceglie la funzione \\u201cR\",\n  \"L\\u2019operatore sceglie la funzione \\u201cVisualizza\",\n];\n\nfor (const searchText of targets) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Target text not found: \" + searchText);\n  }\n\n  const matchRange = results.items[0];\n\n  // Split the match on the apostrophe so we isolate the leading \"L\" as its\n  // own range, leaving \"operatore\"/\"utente\" untouched.\n  const parts = matchRange.split([\"\\u2019\"], false, true, false);\n  parts.load(\"items\");\n  await context.sync();\n\n  const leadingL = parts.items[0];\n  leadingL.insertText(\"<<Include>> Login, dopo cui l\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Add \"<<Include>> Login, dopo cui \" before the actor reference in the\n# first step of each of the three use-case tables, turning:\n#   \"1. L'operatore sceglie ...\"  -> \"1. <<Include>> Login, dopo cui l'operatore sceglie ...\"\n#   \"1. L'utente sceglie ...\"     -> \"1. <<Include>> Login, dopo cui l'utente sceglie ...\"\n#\n# Only the leading capital \"L\" of the actor word is replaced (turned into a\n# lower-case \"l\" and prefixed with the \"<<Include>> Login, dopo cui \" text);\n# the rest of each sentence is left untouched.\n\n$d = $word.ActiveDocument\n\n$apostrophe = [char]0x2019   # '\n$openQuote  = [char]0x201C   # \"\n\n# Unique snippets that identify the first step of each of the three use\n# cases (the text right after the leading \"L\" we need to touch).\n$targets = @(\n  (\"L\" + $apostrophe + \"operatore sceglie la funzione \" + $openQuote + \"Aggiungi\"),\n  (\"L\" + $apostrophe + \"utente sceglie la funzione \"     + $openQuote + \"R\"),\n  (\"L\" + $apostrophe + \"operatore sceglie la funzione \" + $openQuote + \"Visualizza\")\n)\n\nforeach ($targetText in $targets) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Text = $targetText\n    $found = $rng.Find.Execute()\n    if (-not $found) {\n        throw (\"Target text not found: \" + $targetText)\n    }\n\n    # $rng now covers the matched text; isolate just the leading \"L\".\n    $lRange = $d.Range($rng.Start, $rng.Start + 1)\n    if ($lRange.Text -ne \"L\") {\n        throw (\"Expected leading 'L', got: \" + $lRange.Text)\n    }\n\n    $lRange.Text = \"<<Include>> Login, dopo cui l\"\n}\n"}
